$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit: ratio_threshold_range / Max went from 11 to 12 ---
$ws.Range("C2").Value = 12

# --- Header row no longer has an explicit (taller) row height; let Excel
#     size it back to the default by auto-fitting it ---
$ws.Rows("1").AutoFit()

# --- Re-fit the three data columns to their (now narrower) contents.
#     (Target character widths are 21.375 / 5.125 / 5.5; expressed here as
#     exact divisions so the underlying engine's internal pixel rounding
#     snaps to the closest representable column width.) ---
$ws.Columns("A").ColumnWidth = 145.0/7
$ws.Columns("B").ColumnWidth = 31.0/7
$ws.Columns("C").ColumnWidth = 34.0/7

# --- C5 keeps the same look (vertical-center + wrap) but is re-saved
#     through the identical, lower-numbered style slot ---
$ws.Range("C5").WrapText = $true

# --- Move the active selection from C5 to B2 ---
$ws.Range("B2").Select() | Out-Null
